$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet holds weekly price records (one row per market/quality pair) for
# "Apio" sorted by descending date, each new week's data prepended to the top
# of the data block (row 28 onward; rows 2-27 already hold the most recent
# weeks). This commit adds the newest week (2021-09-14) as two new rows
# (Primera / Segunda) and pushes all the previously-existing rows down by two,
# growing the sheet from A1:R120 to A1:R122.

$ws.Rows("28:29").Insert()

# New row 28: Vega Central Mapocho de Santiago, 2021-09-14, Primera
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44453
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100112017
$ws.Range("G28").Value = "Apio"
$ws.Range("H28").Value = "Americana (o)"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 61
$ws.Range("K28").Value = 8000
$ws.Range("L28").Value = 9000
$ws.Range("M28").Value = 8492
$ws.Range("N28").Value = "`$/docena de matas"
$ws.Range("O28").Value = "Región de Coquimbo"
$ws.Range("P28").Value = 1415
$ws.Range("Q28").Value = 6
$ws.Range("R28").Value = "Hortaliza"

# New row 29: Vega Central Mapocho de Santiago, 2021-09-14, Segunda
$ws.Range("A29").Value = 9
$ws.Range("B29").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44453
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = 100112017
$ws.Range("G29").Value = "Apio"
$ws.Range("H29").Value = "Americana (o)"
$ws.Range("I29").Value = "Segunda"
$ws.Range("J29").Value = 43
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 6488
$ws.Range("N29").Value = "`$/docena de matas"
$ws.Range("O29").Value = "Región de Coquimbo"
$ws.Range("P29").Value = 1081
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = "Hortaliza"
